$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.969.96"
$ws.Range("E2").Value = "  -0.14%  "
$ws.Range("D3").Value = "1.955.71"
$ws.Range("E3").Value = "  -0.50%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'243.77"
$ws.Range("E5").Value = "  -1.68%  "
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("D7").Value = "'0.4862"
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "'0.2937"
$ws.Range("E8").Value = "  -0.69%  "
$ws.Range("D9").Value = "'0.07076"
$ws.Range("E9").Value = "  +3.38%  "
$ws.Range("D10").Value = "'19.55"
$ws.Range("E10").Value = "  +1.91%  "
$ws.Range("D11").Value = "'107.60"
$ws.Range("E11").Value = "  +0.49%  "
$ws.Range("D12").Value = "1.956.71"
$ws.Range("E12").Value = "  -0.42%  "
$ws.Range("D13").Value = "'0.07764"
$ws.Range("E13").Value = "  -0.20%  "
$ws.Range("D14").Value = "'5.359"
$ws.Range("E14").Value = "  -1.73%  "
$ws.Range("D15").Value = "'0.7007"
$ws.Range("E15").Value = "  -0.59%  "
$ws.Range("D16").Value = "'277.88"
$ws.Range("E16").Value = "  -3.17%  "
$ws.Range("D17").Value = "30.975.96"
$ws.Range("E17").Value = "  -0.17%  "
$ws.Range("D18").Value = "'0.000007809"
$ws.Range("E18").Value = "  +0.82%  "
$ws.Range("D19").Value = "'13.29"
$ws.Range("E19").Value = "  +0.34%  "
$ws.Range("D20").Value = "2.213.31"
$ws.Range("E21").Value = "  -0.03%  "
$ws.Range("D22").Value = "'5.486"
$ws.Range("E22").Value = "  -2.62%  "
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").Value = "'6.505"
$ws.Range("E24").Value = "  -1.67%  "
$ws.Range("D25").Value = "'9.760"
$ws.Range("E25").Value = "  -2.86%  "
$ws.Range("D26").Value = "'169.09"
$ws.Range("E26").Value = "  -1.00%  "
$ws.Range("D27").Value = "'19.69"
$ws.Range("E27").Value = "  -1.94%  "
$ws.Range("D28").Value = "'2.169"
$ws.Range("E28").Value = "  -1.25%  "
$ws.Range("D29").Value = "'0.1050"
$ws.Range("E29").Value = "  -1.78%  "
$ws.Range("D30").Value = "'1.398"
$ws.Range("E30").Value = "  -3.70%  "
$ws.Range("D31").Value = "'1.565"
$ws.Range("E31").Value = "  -2.14%  "
$ws.Range("D32").Value = "'4.591"
$ws.Range("E32").Value = "  -5.18%  "
$ws.Range("D33").Value = "'4.414"
$ws.Range("E33").Value = "  -2.18%  "
$ws.Range("D34").Value = "'0.04891"
$ws.Range("E34").Value = "  -4.14%  "
$ws.Range("D35").Value = "'0.7529"
$ws.Range("E35").Value = "  -2.84%  "
$ws.Range("D36").Value = "'1.167"
$ws.Range("E36").Value = "  -0.50%  "
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("D38").Value = "'0.02000"
$ws.Range("E38").Value = "  -2.09%  "
$ws.Range("D39").Value = "'2.681"
$ws.Range("E39").Value = "  -1.71%  "
$ws.Range("D40").Value = "'6.541"
$ws.Range("E40").Value = "  +0.29%  "
$ws.Range("D41").Value = "'78.50"
$ws.Range("E41").Value = "  +7.77%  "
$ws.Range("E42").Value = "  -0.65%  "
$ws.Range("D43").Value = "'0.8970"
$ws.Range("E43").Value = "  +0.62%  "
$ws.Range("D44").Value = "'109.34"
$ws.Range("E44").Value = "  -0.60%  "
$ws.Range("D45").Value = "'0.4453"
$ws.Range("E45").Value = "  -0.78%  "
$ws.Range("D46").Value = "'7.843"
$ws.Range("E46").Value = "  +4.00%  "
$ws.Range("D47").Value = "'1.001"
$ws.Range("E47").Value = "  -0.04%  "
$ws.Range("D48").Value = "'986.75"
$ws.Range("E48").Value = "  +1.74%  "
$ws.Range("D49").Value = "'0.1252"
$ws.Range("E49").Value = "  -1.31%  "
$ws.Range("D50").Value = "'9.229"
$ws.Range("E50").Value = "  -2.81%  "
$ws.Range("D51").Value = "'35.94"
$ws.Range("E51").Value = "  -0.46%  "
